$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting rows 57-66 down to 58-67
$ws.Rows.Item(57).Insert()

# Fill in the new row 57 with data.
# Unchanged-from-old-row-57 columns: A,B,C,E,F,G,H,J,O,Q,R
$ws.Range("A57").Value = 1
$ws.Range("B57").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C57").Value = "Arica y Parinacota"
$ws.Range("D57").Value = 44644
$ws.Range("E57").Value = 15
$ws.Range("F57").Value = 100114001
$ws.Range("G57").Value = "Papa"
$ws.Range("H57").Value = "Asterix"
$ws.Range("I57").Value = "1a (cosecha)"
$ws.Range("J57").Value = 1000
$ws.Range("K57").Value = 8500
$ws.Range("L57").Value = 9000
$ws.Range("M57").Value = 8750
$ws.Range("N57").Value = "$/saco 25 kilos"
$ws.Range("O57").Value = "Región de Los Lagos"
$ws.Range("P57").Value = 350
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = "Hortaliza"

# Match the date cell style (s="2") used by column D in other rows
$ws.Range("D57").NumberFormat = $ws.Range("D58").NumberFormat
